$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new daily rows of data (dates are serial numbers; column A already
# carries the date number format via its existing style).
$data = @(
    @(54, 46015,   0,   0,  0),
    @(55, 46016,   7,   7,  0),
    @(56, 46017, 515, 496, 19),
    @(57, 46018,   3,   1,  2),
    @(58, 46019,   1,   1,  0),
    @(59, 46020, 546, 501, 45),
    @(60, 46021, 505, 496,  9),
    @(61, 46022,   0,   0,  0),
    @(62, 46023,   1,   1,  0),
    @(63, 46024, 510, 486, 24),
    @(64, 46025,   1,   1,  0),
    @(65, 46026,   0,   0,  0),
    @(66, 46027, 620, 602, 18),
    @(67, 46028, 562, 538, 24)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Update the view: scroll down so row 53 is at the top, and move the active
# selection to N60 (a single cell rather than the old A53:D53 block).
$win = $excel.ActiveWindow
$win.ScrollRow = 53
$win.ScrollColumn = 1
$ws.Range("N60").Select()
